$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 437.3
$ws.Range("C3").Value = 439.2
$ws.Range("C4").Value = 455.2
$ws.Range("C5").Value = 448.3
$ws.Range("C11").Value = 527.2
$ws.Range("C12").Value = 563.1
$ws.Range("C14").Value = 484.7
$ws.Range("C15").Value = 448.3
$ws.Range("C16").Value = 434.3
$ws.Range("C17").Value = 445.9
$ws.Range("C23").Value = 352.9
